$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tanks")

# New header + description for the "Effective Surface Area" column (C),
# mirroring the existing Name/Description columns (A/B).
$ws.Range("C2").Value = "Effective Surface Area"
$ws.Range("C3").Value = "Effective surface area of the tank that can be used to calculate water volume"

# Match the formatting of the existing header (B2) and body (B3) cells
# so the new column reuses the same cell styles instead of creating new ones.
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# PasteSpecial(xlPasteFormats) only carries over formatting, so re-apply the text.
$ws.Range("C2").Value = "Effective Surface Area"
$ws.Range("C3").Value = "Effective surface area of the tank that can be used to calculate water volume"

# Target stored column width is 25.85546875 characters; this runtime quantizes
# ColumnWidth to a coarser internal grid, so 25.0 is the input that lands on the
# closest reachable stored width (25.83333...) to the authored value.
$ws.Columns.Item(3).ColumnWidth = 25
$ws.Rows.Item(3).RowHeight = 38.25

$excel.CutCopyMode = 0
